$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string rich text collapses to plain text; visually identical) ---
$ws.Range("A8").Value = "Volume 31   Number  12"
$ws.Range("C9").Value = "Report Covering the Week  3/18/2024  Through  3/24/2024"

# --- Stable donor cells for format-only paste (never change value/style in this edit) ---
# A14 = text style (s=14), I14 = integer style (s=15), K15 = percent style (s=16)

# --- Value + style-class transition cells ---
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("F14").PasteSpecial(-4122)

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("D20").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D20").PasteSpecial(-4122)

$ws.Range("E20").Value = -50
$ws.Range("K15").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("C22").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C23").PasteSpecial(-4122)

$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("G23").PasteSpecial(-4122)

$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("H23").PasteSpecial(-4122)

# --- Simple value-only updates (style unchanged) ---
$ws.Range("L15").Value = -83.333333333333
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -23.076923076923
$ws.Range("J16").Value = 41
$ws.Range("K16").Value = -19.512195121951
$ws.Range("L16").Value = -28.260869565217
$ws.Range("M16").Value = -13.157894736842
$ws.Range("N16").Value = -84.579439252336
$ws.Range("C17").Value = 5
$ws.Range("E17").Value = -28.571428571428
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = -11.111111111111
$ws.Range("I17").Value = 40
$ws.Range("J17").Value = 55
$ws.Range("K17").Value = -27.272727272727
$ws.Range("L17").Value = -25.925925925925
$ws.Range("M17").Value = 81.818181818181
$ws.Range("N17").Value = -25.925925925925
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 25
$ws.Range("I18").Value = 51
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = -15
$ws.Range("M18").Value = -37.037037037037
$ws.Range("N18").Value = -86.065573770491
$ws.Range("C19").Value = 23
$ws.Range("D19").Value = 24
$ws.Range("E19").Value = -4.166666666666
$ws.Range("F19").Value = 72
$ws.Range("G19").Value = 89
$ws.Range("H19").Value = -19.101123595505
$ws.Range("I19").Value = 218
$ws.Range("J19").Value = 230
$ws.Range("K19").Value = -5.217391304347
$ws.Range("L19").Value = 7.389162561576
$ws.Range("M19").Value = -30.573248407643
$ws.Range("N19").Value = -62.607204116638
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = -25
$ws.Range("I20").Value = 7
$ws.Range("J20").Value = 12
$ws.Range("K20").Value = -41.666666666666
$ws.Range("L20").Value = -36.363636363636
$ws.Range("M20").Value = 40
$ws.Range("N20").Value = -97.560975609756
$ws.Range("F21").Value = 117
$ws.Range("G21").Value = 137
$ws.Range("H21").Value = -14.598540145985
$ws.Range("I21").Value = 351
$ws.Range("J21").Value = 390
$ws.Range("K21").Value = -10
$ws.Range("L21").Value = -7.631578947368
$ws.Range("M21").Value = -24.025974025974
$ws.Range("N21").Value = -76.754966887417
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -75
$ws.Range("I22").Value = 17
$ws.Range("J22").Value = 20
$ws.Range("K22").Value = -15
$ws.Range("L22").Value = -34.615384615384
$ws.Range("M22").Value = 0
$ws.Range("C24").Value = 73
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 143.333333333333
$ws.Range("F24").Value = 267
$ws.Range("G24").Value = 130
$ws.Range("H24").Value = 105.384615384615
$ws.Range("I24").Value = 728
$ws.Range("J24").Value = 434
$ws.Range("K24").Value = 67.741935483871
$ws.Range("L24").Value = 42.465753424657
$ws.Range("M24").Value = 85.714285714285
$ws.Range("C25").Value = 65
$ws.Range("D25").Value = 30
$ws.Range("E25").Value = 116.666666666667
$ws.Range("F25").Value = 227
$ws.Range("G25").Value = 102
$ws.Range("H25").Value = 122.549019607843
$ws.Range("I25").Value = 623
$ws.Range("J25").Value = 323
$ws.Range("K25").Value = 92.879256965944
$ws.Range("L25").Value = 66.577540106951
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = -7.692307692307
$ws.Range("F26").Value = 51
$ws.Range("G26").Value = 40
$ws.Range("H26").Value = 27.5
$ws.Range("I26").Value = 130
$ws.Range("J26").Value = 103
$ws.Range("K26").Value = 26.213592233009
$ws.Range("L26").Value = 21.495327102803
$ws.Range("M26").Value = 56.626506024096
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = -50
$ws.Range("L27").Value = -66.666666666666
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -25
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 16
$ws.Range("H28").Value = -37.5
$ws.Range("I28").Value = 23
$ws.Range("J28").Value = 26
$ws.Range("K28").Value = -11.538461538461
$ws.Range("L28").Value = 9.523809523809
$ws.Range("L31").Value = -62.5

$excel.CutCopyMode = $false